# Applies the crypto price/volume updates described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    # Force the cell to stay a text value even when it looks like a plain number
    # (Excel would otherwise silently convert strings such as "207.14" to a number).
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $val
}

$ws.Range("D2").Value = "27.660.54"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").Value = "1.587.78"
$ws.Range("E4").Value = "  +0.35%  "
Set-TextCell "D5" "207.14"
$ws.Range("E5").Value = "  -2.05%  "
$ws.Range("E6").Value = "  -3.61%  "
$ws.Range("E7").Value = "  +0.41%  "
Set-TextCell "D8" "22.23"
$ws.Range("E8").Value = "  -4.53%  "
Set-TextCell "D9" "0.252"
$ws.Range("E9").Value = "  -2.11%  "
$ws.Range("E10").Value = "  -2.83%  "
Set-TextCell "D11" "0.0867"
$ws.Range("E11").Value = "  -1.59%  "
$ws.Range("D12").Value = "1.813.06"
$ws.Range("E12").Value = "  -2.61%  "
$ws.Range("D13").Value = "1.592.64"
$ws.Range("E13").Value = "  -2.23%  "
$ws.Range("E14").Value = "  -4.02%  "
$ws.Range("E15").Value = "  -4.74%  "
$ws.Range("D16").Value = "27.652.08"
$ws.Range("E16").Value = "  -0.89%  "
Set-TextCell "D17" "63.44"
$ws.Range("E17").Value = "  -2.34%  "
Set-TextCell "D18" "219.44"
$ws.Range("E18").Value = "  -4.00%  "
$ws.Range("E19").Value = "  -3.29%  "
Set-TextCell "D20" "7.31"
$ws.Range("E20").Value = "  -4.08%  "
$ws.Range("E21").Value = "  +0.45%  "
$ws.Range("E22").Value = "  -4.83%  "
Set-TextCell "D23" "9.65"
$ws.Range("E23").Value = "  -2.83%  "
$ws.Range("E24").Value = "  -3.69%  "
Set-TextCell "D25" "153.62"
$ws.Range("E25").Value = "  -1.08%  "
$ws.Range("E26").Value = "  -1.54%  "
$ws.Range("E27").Value = "  +0.40%  "
Set-TextCell "D28" "15.12"
$ws.Range("E28").Value = "  -2.16%  "
$ws.Range("E29").Value = "  -4.83%  "
$ws.Range("E30").Value = "  -2.73%  "
$ws.Range("E31").Value = "  -2.78%  "
Set-TextCell "D32" "3.23"
$ws.Range("E32").Value = "  -5.41%  "
$ws.Range("D33").Value = "1.369.30"
$ws.Range("E33").Value = "  -3.42%  "
$ws.Range("E34").Value = "  -5.59%  "
Set-TextCell "D35" "1.54"
$ws.Range("E35").Value = "  -4.70%  "
$ws.Range("E36").Value = "  -3.45%  "
$ws.Range("E37").Value = "  -0.61%  "
Set-TextCell "D38" "0.0167"
$ws.Range("E38").Value = "  -1.30%  "
Set-TextCell "D39" "0.535"
$ws.Range("E39").Value = "  -3.26%  "
Set-TextCell "D40" "0.825"
$ws.Range("E40").Value = "  -3.37%  "
$ws.Range("E41").Value = "  +0.40%  "
Set-TextCell "D42" "0.972"
Set-TextCell "D43" "64.15"
$ws.Range("E43").Value = "  -2.77%  "
$ws.Range("E44").Value = "  +2.13%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell "D45" "5.18"
$ws.Range("E45").Value = "  -4.43%  "
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "1.724.33"
$ws.Range("E46").Value = "  -2.62%  "
Set-TextCell "D47" "1.72"
$ws.Range("E47").Value = "  -5.02%  "
Set-TextCell "D48" "87.38"
$ws.Range("E48").Value = "  -1.48%  "
$ws.Range("D49").Value = "0.0₆0100"
$ws.Range("E49").Value = "  -1.68%  "
Set-TextCell "D50" "0.0965"
$ws.Range("E50").Value = "  -4.58%  "
Set-TextCell "D51" "0.0494"
$ws.Range("E51").Value = "  -1.57%  "
